$d = $word.ActiveDocument

# Locate the paragraph that holds the "LOT2053: Microbiologia (Requisito fraco)"
# line (use Find to make sure the text actually exists in the document before
# we touch anything).
$anchor = $d.Content
$anchorFound = $anchor.Find.Execute("LOT2053: Microbiologia (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($anchorFound) {
    # Walk the Paragraphs collection to find the 1-based index of that
    # paragraph; the three paragraphs right after it (a blank paragraph,
    # "Ver no Jupiter Salvar em pdf Salvar em docx" and the
    # "© 2020 . Contact: ..." footer line) are the ones being removed.
    $count = $d.Paragraphs.Count
    $lot2053Index = -1
    for ($i = 1; $i -le $count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.Contains("LOT2053: Microbiologia (Requisito fraco)")) {
            $lot2053Index = $i
            break
        }
    }

    if ($lot2053Index -gt 0) {
        $delStart = $lot2053Index + 1
        $delEnd = $lot2053Index + 3

        $rangeStart = $d.Paragraphs.Item($delStart).Range.Start
        $rangeEnd = $d.Paragraphs.Item($delEnd).Range.End

        $killRange = $d.Range($rangeStart, $rangeEnd)
        $killRange.Delete()
    }
}
